$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows for "103_TruckInsurance_002_VehicleData_*" testcases
#     (these land as rows 15-16 after the SendQuote rows are inserted above them) ---
$ws.Rows("13:14").Insert()
$ws.Cells.Item(13, 1).Value2 = "103_TruckInsurance_002_VehicleData_001_MandatoryFields"
$ws.Cells.Item(14, 1).Value2 = "103_TruckInsurance_002_VehicleData_002_FieldHintsAndErrors"
$ws.Cells.Item(13, 2).Value2 = "var103_TruckInsurance_002_VehicleData_001_MandatoryFields"
$ws.Cells.Item(14, 2).Value2 = "var103_TruckInsurance_002_VehicleData_002_FieldHintsAndErrors"
$ws.Cells.Item(13, 3).Value2 = "Open Truck Insurance"
$ws.Cells.Item(14, 3).Value2 = "Open Truck Insurance"
$ws.Cells.Item(13, 5).Value2 = "103_TruckInsurance_002_VehicleData_001_MandatoryFields"
$ws.Cells.Item(14, 5).Value2 = "103_TruckInsurance_002_VehicleData_002_FieldHintsAndErrors"

# --- Insert 2 new rows for "102_AutomobileInsurance_006_SendQuote_*" testcases ---
$ws.Rows("12:13").Insert()
$ws.Cells.Item(12, 1).Value2 = "102_AutomobileInsurance_006_SendQuote_001_MandatoryFields"
$ws.Cells.Item(13, 1).Value2 = "102_AutomobileInsurance_006_SendQuote_002_FieldHintsAndErrors"
$ws.Cells.Item(12, 2).Value2 = "var102_AutomobileInsurance_006_SendQuote_001_MandatoryFields"
$ws.Cells.Item(13, 2).Value2 = "var102_AutomobileInsurance_006_SendQuote_002_FieldHintsAndErrors"
$ws.Cells.Item(12, 3).Value2 = "Open Automobile Insurance"
$ws.Cells.Item(13, 3).Value2 = "Open Automobile Insurance"
$ws.Cells.Item(12, 4).Value2 = "102_AutomobileInsurance_006_SendQuote_001_MandatoryFields"
$ws.Cells.Item(13, 4).Value2 = "102_AutomobileInsurance_006_SendQuote_002_FieldHintsAndErrors"

# --- Column E grew wider to fit the new longer testcase names ---
$ws.Columns.Item(5).ColumnWidth = 53

# --- Update the saved selection/active cell ---
[void]$ws.Range("D27").Select()
